$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.937.32'
$ws.Range('E2').Value = '  +1.24%  '

$ws.Range('D3').Value = '1.708.21'
$ws.Range('E3').Value = '  +0.90%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4024'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4081'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.97%  '

$ws.Range('B9').Value = 'Polygon'
$ws.Range('C9').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.483'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.65%  '

$ws.Range('B10').Value = 'BinanceUSD'
$ws.Range('C10').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.003'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.09%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.83'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.42%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08842'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.20%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.36'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.503'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.77%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.042'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001346'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.50%  '

$ws.Range('D17').Value = '1.661.81'
$ws.Range('E17').Value = '  -1.73%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.46'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07182'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.93%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.281'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.31%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.005'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.93%  '

$ws.Range('D24').Value = '24.946.47'
$ws.Range('E24').Value = '  +1.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.338'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.58%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.906'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.39%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.432'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +22.74%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.96%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.84%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '143.77'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.97%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.216'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.67%  '

$ws.Range('E32').Value = '  +14.11%  '

$ws.Range('D33').Value = '1.847.43'
$ws.Range('E33').Value = '  -1.56%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08752'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.20%  '

$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.374'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.07%  '

$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.03200'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.47%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.032'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.93%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2879'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.40%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8503'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.16%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.88'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09465'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.94%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.91%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.480'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.99%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.76'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.18%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.731'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.43%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7483'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.241'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.85%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.394'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.90%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.10%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.03%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.08415'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.57%  '
